# Update the ABUNDANCE column (column 2) values in the single data table
# per the new correlation results referenced in the commit message.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of row index -> [old value, new value] for the ABUNDANCE column (col 2)
$updates = @(
    @{ Row = 2;  Old = "0.600"; New = "0.584" },  # CODNEAR
    @{ Row = 3;  Old = "0.339"; New = "0.347" },  # CODNEARNCW
    @{ Row = 4;  Old = "0.348"; New = "0.174" },  # CODFAPL
    @{ Row = 5;  Old = "0.583"; New = "0.490" },  # CODICE
    @{ Row = 6;  Old = "0.199"; New = "0.139" },  # CODBA2532
    @{ Row = 7;  Old = "0.408"; New = "0.520" },  # CODKAT
    @{ Row = 8;  Old = "0.092"; New = "0.138" },  # CODIS
    @{ Row = 9;  Old = "0.204"; New = "0.288" },  # CODVIa
    @{ Row = 10; Old = "0.329"; New = "0.331" },  # CODIIIaW
    @{ Row = 11; Old = "0.764"; New = "0.709" },  # HAKENRTN
    @{ Row = 12; Old = "0.427"; New = "0.140" }   # HAKESOTH
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, 2)
    $rng = $cell.Range
    # Exclude the trailing cell-mark/paragraph-mark from the range so we
    # only touch the numeric text itself.
    $rng.End = $rng.End - 1

    $current = $rng.Text
    if ($current -ne $u.Old) {
        Write-Output "WARNING: row $($u.Row) expected '$($u.Old)' but found '$current'"
    }

    # Assign the replacement text directly on the (already cell-scoped)
    # range. This keeps the edit confined to this single cell, unlike
    # Find.Execute with ReplaceAll + wdFindContinue, which can wrap
    # around and replace matching text elsewhere in the document.
    $rng.Text = $u.New
}

Write-Output "Done updating ABUNDANCE column."
